# Scheduled market-data refresh: updates currentAveragePrice(NQ/HQ) and the
# derived LevePrice/LeveProfit columns (H,I,J,K,L,M,N) across the per-job
# leve tables (one table per crafting job worksheet).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H92").Value = 752.44446
$ws.Range("J92").Value = 1000
$ws.Range("L92").Value = 1000
$ws.Range("N92").Value = -3496

$ws.Range("H96").Value = 1047.909
$ws.Range("I96").Value = 286.7143
$ws.Range("K96").Value = 860.1428999999999
$ws.Range("M96").Value = 512.8571000000001

$ws.Range("H98").Value = 2292.4546
$ws.Range("I98").Value = 1681.2778
$ws.Range("K98").Value = 1681.2778
$ws.Range("M98").Value = -183.2778000000001

$ws.Range("H100").Value = 2111
$ws.Range("I100").Value = 222
$ws.Range("K100").Value = 222
$ws.Range("M100").Value = 319

$ws.Range("H101").Value = 1823.6111
$ws.Range("I101").Value = 1288.0834
$ws.Range("K101").Value = 3864.2502
$ws.Range("M101").Value = -2242.2502

$ws.Range("H106").Value = 5875
$ws.Range("I106").Value = 6348.591
$ws.Range("K106").Value = 6348.591
$ws.Range("M106").Value = -5717.591

$ws.Range("H112").Value = 6137.5854
$ws.Range("J112").Value = 6253.525
$ws.Range("L112").Value = 18760.575
$ws.Range("N112").Value = -20976.575

$ws.Range("H122").Value = 2292.4546
$ws.Range("I122").Value = 1681.2778
$ws.Range("K122").Value = 5043.8334
$ws.Range("M122").Value = -2593.8334

$ws.Range("H138").Value = 3460.4744
$ws.Range("J138").Value = 3644.1
$ws.Range("L138").Value = 10932.3
$ws.Range("N138").Value = -21212.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2392.0952
$ws.Range("I2").Value = 2554.8235
$ws.Range("K2").Value = 2554.8235
$ws.Range("M2").Value = -2441.8235

$ws.Range("H28").Value = 10789.5
$ws.Range("I28").Value = 11343.125
$ws.Range("J28").Value = 8575
$ws.Range("K28").Value = 11343.125
$ws.Range("L28").Value = 8575
$ws.Range("M28").Value = -11151.125
$ws.Range("N28").Value = -8959

$ws.Range("H32").Value = 8907.931
$ws.Range("I32").Value = 5439.831
$ws.Range("K32").Value = 5439.831
$ws.Range("M32").Value = -5152.831

$ws.Range("H45").Value = 6799.3335
$ws.Range("I45").Value = 17028.846
$ws.Range("J45").Value = 3776.9773
$ws.Range("K45").Value = 17028.846
$ws.Range("L45").Value = 3776.9773
$ws.Range("M45").Value = -16651.846
$ws.Range("N45").Value = -4530.9773

$ws.Range("H61").Value = 3223.9167
$ws.Range("I61").Value = 2285
$ws.Range("K61").Value = 2285
$ws.Range("M61").Value = -2073

$ws.Range("H99").Value = 10789.5
$ws.Range("I99").Value = 11343.125
$ws.Range("J99").Value = 8575
$ws.Range("K99").Value = 11343.125
$ws.Range("L99").Value = 8575
$ws.Range("M99").Value = -8348.125
$ws.Range("N99").Value = -14565

$ws.Range("H110").Value = 3812.4211
$ws.Range("I110").Value = 3709.0667
$ws.Range("J110").Value = 4200
$ws.Range("K110").Value = 3709.0667
$ws.Range("L110").Value = 4200
$ws.Range("M110").Value = -1664.0667
$ws.Range("N110").Value = -8290

$ws.Range("H116").Value = 2392.0952
$ws.Range("I116").Value = 2554.8235
$ws.Range("K116").Value = 2554.8235
$ws.Range("M116").Value = -260.8235

$ws.Range("H132").Value = 2457.6667
$ws.Range("I132").Value = 1627
$ws.Range("K132").Value = 4881
$ws.Range("M132").Value = -2351

$ws.Range("H136").Value = 3223.9167
$ws.Range("I136").Value = 2285
$ws.Range("K136").Value = 6855
$ws.Range("M136").Value = -4305

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2392.0952
$ws.Range("I3").Value = 2554.8235
$ws.Range("K3").Value = 2554.8235
$ws.Range("M3").Value = -2440.8235

$ws.Range("H58").Value = 35572.5
$ws.Range("J58").Value = 35572.5
$ws.Range("L58").Value = 35572.5
$ws.Range("N58").Value = -36160.5

$ws.Range("H134").Value = 2630.2222
$ws.Range("I134").Value = 2572.375
$ws.Range("J134").Value = 3093
$ws.Range("K134").Value = 7717.125
$ws.Range("L134").Value = 9279
$ws.Range("M134").Value = -5182.125
$ws.Range("N134").Value = -14349

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1228.6666
$ws.Range("I22").Value = 1424.4
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 1424.4
$ws.Range("L22").Value = 250
$ws.Range("M22").Value = -1074.4
$ws.Range("N22").Value = -950

$ws.Range("H31").Value = 7318.018
$ws.Range("I31").Value = 5802.185
$ws.Range("J31").Value = 8779.714
$ws.Range("K31").Value = 5802.185
$ws.Range("L31").Value = 8779.714
$ws.Range("M31").Value = -5507.185
$ws.Range("N31").Value = -9369.714

$ws.Range("H34").Value = 7318.018
$ws.Range("I34").Value = 5802.185
$ws.Range("J34").Value = 8779.714
$ws.Range("K34").Value = 5802.185
$ws.Range("L34").Value = 8779.714
$ws.Range("M34").Value = -5600.185
$ws.Range("N34").Value = -9183.714

$ws.Range("H62").Value = 2882.3333
$ws.Range("I62").Value = 2882.3333
$ws.Range("K62").Value = 2882.3333
$ws.Range("M62").Value = -2258.3333

$ws.Range("H65").Value = 2882.3333
$ws.Range("I65").Value = 2882.3333
$ws.Range("K65").Value = 14411.6665
$ws.Range("M65").Value = -11291.6665

$ws.Range("H86").Value = 16678435
$ws.Range("I86").Value = 30315704
$ws.Range("K86").Value = 30315704
$ws.Range("M86").Value = -30314581

$ws.Range("H89").Value = 16678435
$ws.Range("I89").Value = 30315704
$ws.Range("K89").Value = 151578520
$ws.Range("M89").Value = -151572904

$ws.Range("H99").Value = 5563957.5
$ws.Range("I99").Value = 13902145
$ws.Range("K99").Value = 13902145
$ws.Range("M99").Value = -13900647

$ws.Range("H126").Value = 5563957.5
$ws.Range("I126").Value = 13902145
$ws.Range("K126").Value = 41706435
$ws.Range("M126").Value = -41703965

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1834
$ws.Range("J68").Value = 1834
$ws.Range("L68").Value = 5502
$ws.Range("N68").Value = -7124

$ws.Range("H71").Value = 1834
$ws.Range("J71").Value = 1834
$ws.Range("L71").Value = 16506
$ws.Range("N71").Value = -24618

$ws.Range("H128").Value = 109979.2
$ws.Range("I128").Value = 109979.2
$ws.Range("K128").Value = 329937.6
$ws.Range("M128").Value = -324957.6

$ws.Range("H132").Value = 2668.6667
$ws.Range("I132").Value = 2499
$ws.Range("J132").Value = 2680.7856
$ws.Range("K132").Value = 22491
$ws.Range("L132").Value = 24127.0704
$ws.Range("M132").Value = -19961
$ws.Range("N132").Value = -29187.0704

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 200.90909
$ws.Range("I2").Value = 244
$ws.Range("K2").Value = 244
$ws.Range("M2").Value = -131

$ws.Range("H113").Value = 13012.556
$ws.Range("J113").Value = 24227.5
$ws.Range("L113").Value = 24227.5
$ws.Range("N113").Value = -28567.5

$ws.Range("H122").Value = 3715.8462
$ws.Range("I122").Value = 2641.9
$ws.Range("J122").Value = 7295.6665
$ws.Range("K122").Value = 7925.700000000001
$ws.Range("L122").Value = 21886.9995
$ws.Range("M122").Value = -5475.700000000001
$ws.Range("N122").Value = -26786.9995

$ws.Range("H126").Value = 4281.5
$ws.Range("I126").Value = 3747.5
$ws.Range("J126").Value = 5349.5
$ws.Range("K126").Value = 11242.5
$ws.Range("L126").Value = 16048.5
$ws.Range("M126").Value = -8772.5
$ws.Range("N126").Value = -20988.5

$ws.Range("H132").Value = 4128.685
$ws.Range("I132").Value = 3801.2559
$ws.Range("J132").Value = 5408.636
$ws.Range("K132").Value = 11403.7677
$ws.Range("L132").Value = 16225.908
$ws.Range("M132").Value = -8873.7677
$ws.Range("N132").Value = -21285.908

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1707.4445
$ws.Range("I22").Value = 1546
$ws.Range("J22").Value = 2999
$ws.Range("K22").Value = 1546
$ws.Range("L22").Value = 2999
$ws.Range("M22").Value = -1251
$ws.Range("N22").Value = -3589

$ws.Range("H27").Value = 1707.4445
$ws.Range("I27").Value = 1546
$ws.Range("J27").Value = 2999
$ws.Range("K27").Value = 1546
$ws.Range("L27").Value = 2999
$ws.Range("M27").Value = -1439
$ws.Range("N27").Value = -3213

$ws.Range("H100").Value = 3151.818
$ws.Range("J100").Value = 1320
$ws.Range("L100").Value = 1320
$ws.Range("N100").Value = -2402

$ws.Range("H132").Value = 23048.938
$ws.Range("I132").Value = 12616.083
$ws.Range("K132").Value = 37848.249
$ws.Range("M132").Value = -35318.249

$ws.Range("H136").Value = 2213
$ws.Range("I136").Value = 2261.077
$ws.Range("J136").Value = 1796.3334
$ws.Range("K136").Value = 6783.231000000001
$ws.Range("L136").Value = 5389.0002
$ws.Range("M136").Value = -4233.231000000001
$ws.Range("N136").Value = -10489.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4870.9614
$ws.Range("I122").Value = 4416.1904
$ws.Range("K122").Value = 13248.5712
$ws.Range("M122").Value = -10798.5712

$ws.Range("H132").Value = 5471.269
$ws.Range("I132").Value = 5182.5
$ws.Range("J132").Value = 6433.8335
$ws.Range("K132").Value = 15547.5
$ws.Range("L132").Value = 19301.5005
$ws.Range("M132").Value = -13017.5
$ws.Range("N132").Value = -24361.5005

$ws.Range("H136").Value = 14439.363
$ws.Range("I136").Value = 20226.8
$ws.Range("K136").Value = 60680.39999999999
$ws.Range("M136").Value = -58130.39999999999
